$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header row (bold, centered, bordered - same as H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new I and J column data for rows 2-34.
$data = @(
    @(2,2,6),
    @(3,6,7),
    @(4,3,7),
    @(5,5,8),
    @(6,1,3),
    @(7,3,6),
    @(8,6,9),
    @(9,7,7),
    @(10,2,3),
    @(11,9,9),
    @(12,5,9),
    @(13,4,7),
    @(14,4,6),
    @(15,5,7),
    @(16,1,5),
    @(17,1,4),
    @(18,3,7),
    @(19,3,4),
    @(20,2,7),
    @(21,1,4),
    @(22,1,5),
    @(23,1,5),
    @(24,1,6),
    @(25,1,5),
    @(26,1,5),
    @(27,7,7),
    @(28,1,3),
    @(29,1,5),
    @(30,1,4),
    @(31,1,6),
    @(32,1,5),
    @(33,1,3),
    @(34,1,2)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 9).Value = $entry[1]
    $ws.Cells.Item($row, 10).Value = $entry[2]
}

Write-Output "Added columns I and J (I0 / IF) to sheet1"
